# Apply cost-data adaptation to the embodied emissions systems workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New combined source-cost citation used for the heat-generation rows (2-9)
$newConversionSource = "T. Gürber (2020), D. Sigrist (2018), WWF (2019), Energieschweiz und erneuerbarheizen (2020), energie 360°"
# New short citation used for rows 13, 14 and 17
$shortSource = "T. Gürber (2020)"

# --- Row 2: Pellets ---
$ws.Range("K2").Value = 2110
$ws.Range("L2").Value = $newConversionSource

# --- Row 3: Wood ---
$ws.Range("K3").Value = 2110
$ws.Range("L3").Value = $newConversionSource

# --- Row 4: ASHP ---
$ws.Range("K4").Value = 2290
$ws.Range("L4").Value = $newConversionSource

# --- Row 5: GSHP ---
$ws.Range("K5").Value = 3970
$ws.Range("L5").Value = $newConversionSource

# --- Row 6: electric (cost unchanged, only source text changes) ---
$ws.Range("L6").Value = $newConversionSource

# --- Row 7: Natural Gas ---
$ws.Range("K7").Value = 1220
$ws.Range("L7").Value = $newConversionSource

# --- Row 8: Oil ---
$ws.Range("K8").Value = 1040
$ws.Range("L8").Value = $newConversionSource

# --- Row 9: district (cost unchanged, only source text changes) ---
$ws.Range("L9").Value = $newConversionSource

# --- Row 13: floor heating (source citation shortened) ---
$ws.Range("L13").Value = $shortSource

# --- Row 14: radiator (source citation shortened) ---
$ws.Range("L14").Value = $shortSource

# --- Row 17: mechanical ventilation (source citation shortened) ---
$ws.Range("L17").Value = $shortSource

# --- Sheet view adjustments ---
$ws.Range("K22").Select()
